# Workshop documentation and typo fix.
#
# "Tripple Residential Pod" -> "Triple Residential Pod" everywhere it
# appears in the workbook (the Property Types lookup table, and every
# Properties-table row that used that property type). Using a
# workbook-wide Find & Replace mirrors exactly what happened upstream:
# the shared string gets corrected in place for every cell that held it.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Tripple Residential Pod", "Triple Residential Pod")
}

# --- restore the cursor/selection left on each sheet ---
$wsLeases = $wb.Worksheets.Item("Leases")
$wsProperties = $wb.Worksheets.Item("Properties")
$wsPropertyTypes = $wb.Worksheets.Item("Property Types")

$wsPropertyTypes.Activate()
$wsPropertyTypes.Range("B11").Select()

$wsProperties.Activate()
$wsProperties.Range("A2").Select()

$wsLeases.Activate()
$wsLeases.Range("A6").Select()
